# Regenerate the s_vals data to filter save games.
# Updates the numeric values in B2:G4 on the active sheet to the newly
# computed figures, leaving headers (row 1) and dates (column A) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B2" = 0.6606524410359556
    "C2" = 1.655778082260271
    "D2" = 0.1494219747398047
    "E2" = 0.4942365360607697
    "F2" = 1
    "G2" = 2.960089034096801

    "B3" = 3.286832544864788
    "C3" = 1.655778082260271
    "D3" = 6708.013860684405
    "E3" = 10.19245300693656
    "F3" = 0
    "G3" = 6723.148924318466

    "B4" = 3.286832544864788
    "C4" = 117.745847958593
    "D4" = 0.7527432677738641
    "E4" = 10.19245300693656
    "F4" = 1
    "G4" = 131.9778767781682
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
